$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
# existing "2022-Q3" sheet) - mirrors the workbook.xml <sheets> reorder.
#
# NOTE: sheet object handles obtained *before* a Worksheets.Add/Move
# re-resolve by position, not by stable identity - after the collection
# is mutated they can silently refer to a different sheet. So the add
# is done first, in isolation, and every sheet handle used afterwards is
# re-fetched by name.
# ---------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($insertBefore)
$q4Sheet.Name = "2022-Q4"

$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2022-Q4")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# Populate "2022-Q4" with the same layout/styling as the "2022-Q3" sheet
# (header row bold+bordered, data rows plain).
# ---------------------------------------------------------------------
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

function Set-FundRow {
    param($ws, $row, $aStyleSrc, $code, $fundName, $scale, $totalPos, $posPct, $value, $rank)

    # A column: index number, copy the numeric+style formatting from the
    # reference sheet so it reuses the existing cellXf (s="2") instead of
    # minting a new style.
    $ws.Range("A" + $row).Value = $row - 2
    $aStyleSrc.Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)  # xlPasteFormats

    # B holds a numeric-looking fund code (e.g. "015641") that must stay
    # TEXT (with the leading zero) and unstyled - same trick as the D/E/F/G
    # columns below.
    $ws.Range("B" + $row).Value = "'" + $code
    $ws.Range("B" + $row).Style = "Normal"

    # C is the fund name - plain Chinese text, not numeric-looking, so a
    # normal assignment already keeps it as an unstyled string.
    $ws.Range("C" + $row).Value = $fundName

    # D/E/F/G hold numeric-looking values but must stay TEXT (matches the
    # inlineStr cells in the source file) with no style override. Force
    # text via a leading apostrophe, then reset the style back to Normal
    # so no extra Text-format cellXf gets minted.
    $ws.Range("D" + $row).Value = "'" + $scale
    $ws.Range("D" + $row).Style = "Normal"
    $ws.Range("E" + $row).Value = "'" + $totalPos
    $ws.Range("E" + $row).Style = "Normal"
    $ws.Range("F" + $row).Value = "'" + $posPct
    $ws.Range("F" + $row).Style = "Normal"
    $ws.Range("G" + $row).Value = "'" + $value
    $ws.Range("G" + $row).Style = "Normal"

    # H is a real number, no style.
    $ws.Range("H" + $row).Value = $rank
}

$aRef = $q3Sheet.Range("A2")
Set-FundRow $q4Sheet 2 $aRef "015641" "银华数字经济股票A"     "0.20" "94.75" "5.08" "0.0102" 10
Set-FundRow $q4Sheet 3 $aRef "002068" "东方多策略灵活配置混合C" "0.26" "87.87" "2.76" "0.0072" 8
Set-FundRow $q4Sheet 4 $aRef "015642" "银华数字经济股票C"     "0.06" "94.75" "5.08" "0.0030" 10
Set-FundRow $q4Sheet 5 $aRef "400023" "东方多策略灵活配置混合A" "0.03" "87.87" "2.76" "0.0008" 8

# ---------------------------------------------------------------------
# Update "总计": existing rows shift down one quarter, a new row is
# appended for 2022-Q2 that mirrors the old 2022-Q2 row.
# ---------------------------------------------------------------------
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 0.02

$totalSheet.Range("B3").Value = "2022-Q3"

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.01

# Restore the originally active sheet/tab.
$totalSheet.Activate()
